$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct total marks error on the "Marking" and "Total" rows
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 76
$ws.Range("E12").Value = "76 / 112"
